# Auto-generated edit script: updates cryptocurrency Price (D) and
# Volume(1h) (E) columns on Sheet1 to match the refreshed scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.526.35'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '2.497.72'
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''574.90'
$ws.Range("E5").Value = '  -0.71%  '
$ws.Range("D6").Value = '''166.43'
$ws.Range("E6").Value = '  -0.56%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '''0.512'
$ws.Range("E8").Value = '  -1.61%  '
$ws.Range("D9").Value = '2.496.06'
$ws.Range("E9").Value = '  -0.52%  '
$ws.Range("E10").Value = '  +0.07%  '
$ws.Range("E11").Value = '  +0.04%  '
$ws.Range("E12").Value = '  +3.08%  '
$ws.Range("D13").Value = '''4.92'
$ws.Range("E13").Value = '  +1.10%  '
$ws.Range("D14").Value = '2.954.89'
$ws.Range("E14").Value = '  -0.50%  '
$ws.Range("D15").Value = '69.386.11'
$ws.Range("E15").Value = '  -0.19%  '
$ws.Range("E16").Value = '  +1.26%  '
$ws.Range("D17").Value = '''24.69'
$ws.Range("E17").Value = '  -1.12%  '
$ws.Range("D18").Value = '2.518.00'
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("D19").Value = '''11.19'
$ws.Range("E19").Value = '  -1.11%  '
$ws.Range("D20").Value = '''7.44'
$ws.Range("E20").Value = '  -4.87%  '
$ws.Range("D21").Value = '''347.40'
$ws.Range("E21").Value = '  -0.37%  '
$ws.Range("E22").Value = '  -1.07%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").Value = '''70.69'
$ws.Range("E25").Value = '  +2.83%  '
$ws.Range("D26").Value = '''3.94'
$ws.Range("E26").Value = '  -1.02%  '
$ws.Range("D27").Value = '''8.73'
$ws.Range("E27").Value = '  -2.54%  '
$ws.Range("D28").Value = '2.627.90'
$ws.Range("E28").Value = '  -0.45%  '
$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  +0.70%  '
$ws.Range("E30").Value = '  -1.55%  '
$ws.Range("D31").Value = '''7.84'
$ws.Range("E31").Value = '  -0.67%  '
$ws.Range("D32").Value = '''456.37'
$ws.Range("E32").Value = '  -1.92%  '
$ws.Range("E33").Value = '  -4.65%  '
$ws.Range("E34").Value = '  -1.27%  '
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("E36").Value = '  -2.07%  '
$ws.Range("D37").Value = '''156.13'
$ws.Range("E37").Value = '  +1.16%  '
$ws.Range("D38").Value = '''19.04'
$ws.Range("E38").Value = '  +0.34%  '
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("E40").Value = '  +0.03%  '
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("D42").Value = '''4.67'
$ws.Range("E42").Value = '  -1.83%  '
$ws.Range("D43").Value = '''1.59'
$ws.Range("E43").Value = '  -0.67%  '
$ws.Range("D44").Value = '''38.08'
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("E45").Value = '  -4.37%  '
$ws.Range("E46").Value = '  -6.43%  '
$ws.Range("D47").Value = '''141.22'
$ws.Range("E47").Value = '  -1.19%  '
$ws.Range("D48").Value = '''3.48'
$ws.Range("E48").Value = '  -0.33%  '
$ws.Range("E49").Value = '  -2.30%  '
$ws.Range("D50").Value = '''0.0730'
$ws.Range("E50").Value = '  -0.32%  '
$ws.Range("E51").Value = '  -0.31%  '
